# Node size is variable regarding its degree: update the "size" column
# (F2:F33) on Sheet2 from a constant 1 to a constant 10, and leave the
# sheet with F2:F33 selected (mirroring the author's manual edit/selection
# captured in the saved workbook view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Update the F column (size) values for data rows 2-33 from 1 to 10.
$ws.Range("F2:F33").Value = 10

# Reflect the resulting selection/active cell state on the sheet
# (active cell F2, with the full F2:F33 range selected).
$ws.Range("F2:F33").Select()
